$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths for the new "common mistakes" table (J:O)
$ws.Columns.Item(10).ColumnWidth = 14.33203125
$ws.Columns.Item(11).ColumnWidth = 25.1640625
$ws.Columns.Item(12).ColumnWidth = 18.83203125
$ws.Columns.Item(13).ColumnWidth = 16.83203125
$ws.Columns.Item(14).ColumnWidth = 13.83203125
$ws.Columns.Item(15).ColumnWidth = 14

# Header row (row 6) for the new table
$ws.Range("J6").Value = "key_id"
$ws.Range("K6").Value = "roof_type"
$ws.Range("L6").Value = "wall type"
$ws.Range("M6").Value = "floor type"
$ws.Range("N6").Value = "rooms"
$ws.Range("O6").Value = "inc_barn"

# Merged title cell above the table (row 5)
$ws.Range("K5:N5").Merge()

# Data rows (7-16)
$data = @(
  @(1, "grass", "muddaub", "errth", 1, "no"),
  @(2, "grass", "muddaub", "earth", 1, "no"),
  @(3, "mabati_sloping", "burntbricks", "cement", -99, "no"),
  @(4, "mabatisloping", "burntbricks", "earth", 1, "no"),
  @(5, "grass", "burntbricks", "earth", 1, "no"),
  @(6, "grass", "muddaub", "earth", 1, "no"),
  @(7, "grass", "muddaub", "earth", 1, "no"),
  @(8, "mabatisloping", "burntbricks", "cement", 3, "no"),
  @(9, "grass", "burntbricks", "earth", 1, "no"),
  @(10, "mabatisloping", "burntbricks", "cement", 5, "yes")
)

$row = 7
foreach ($rec in $data) {
  $ws.Range("J$row").Value = $rec[0]
  $ws.Range("K$row").Value = $rec[1]
  $ws.Range("L$row").Value = $rec[2]
  $ws.Range("M$row").Value = $rec[3]
  $ws.Range("N$row").Value = $rec[4]
  $ws.Range("O$row").Value = $rec[5]
  $row++
}

$ws.Range("N26").Select()
